$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (in-place character replace, same length) ---
$ws.Range("A8").Characters(21, 2).Text = "17"
$ws.Range("C9").Characters(27, 9).Text = "4/21/2025"
$ws.Range("C9").Characters(47, 9).Text = "4/27/2025"

# --- Style donor NumberFormat strings captured once from stable cells ---
$fmt_C = $ws.Range("C16").NumberFormat
$fmt_D = $ws.Range("D16").NumberFormat
$fmt_E = $ws.Range("E16").NumberFormat
$fmt_F = $ws.Range("F16").NumberFormat
$fmt_G = $ws.Range("G16").NumberFormat
$fmt_H = $ws.Range("H16").NumberFormat

# --- Cell value + style updates ---
# Row 14
$ws.Range("C14").NumberFormat = $fmt_C
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 5
$ws.Range("L14").Value = 400
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = -64.285714285714

# Row 15
$ws.Range("D15").NumberFormat = $fmt_D
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = $fmt_E
$ws.Range("E15").Value = -100
$ws.Range("G15").NumberFormat = $fmt_G
$ws.Range("G15").Value = 1
$ws.Range("H15").NumberFormat = $fmt_H
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -60
$ws.Range("M15").Value = -55.555555555555

# Row 16
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -41.666666666666
$ws.Range("I16").Value = 56
$ws.Range("J16").Value = 118
$ws.Range("K16").Value = -52.542372881355
$ws.Range("L16").Value = -41.052631578947
$ws.Range("M16").Value = -54.098360655737
$ws.Range("N16").Value = -87.719298245614

# Row 17
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -35.483870967741
$ws.Range("I17").Value = 98
$ws.Range("J17").Value = 134
$ws.Range("K17").Value = -26.865671641791
$ws.Range("L17").Value = -20.967741935483
$ws.Range("M17").Value = 4.255319148936
$ws.Range("N17").Value = -64.233576642335

# Row 18
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -45
$ws.Range("I18").Value = 48
$ws.Range("J18").Value = 78
$ws.Range("K18").Value = -38.461538461538
$ws.Range("L18").Value = -34.246575342465
$ws.Range("M18").Value = -49.473684210526
$ws.Range("N18").Value = -88.914549653579

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 44.444444444444
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -22.448979591836
$ws.Range("I19").Value = 156
$ws.Range("J19").Value = 194
$ws.Range("K19").Value = -19.587628865979
$ws.Range("L19").Value = -31.277533039647
$ws.Range("M19").Value = 71.428571428571
$ws.Range("N19").Value = -4.878048780487

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 7.692307692307
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 62
$ws.Range("K20").Value = -24.193548387096
$ws.Range("L20").Value = -11.320754716981
$ws.Range("M20").Value = 2.173913043478
$ws.Range("N20").Value = -85.448916408668

# Row 21
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -6.25
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = -28.260869565217
$ws.Range("I21").Value = 414
$ws.Range("J21").Value = 596
$ws.Range("K21").Value = -30.536912751677
$ws.Range("L21").Value = -28.373702422145
$ws.Range("M21").Value = -9.803921568627
$ws.Range("N21").Value = -75.459395376407

# Row 22
$ws.Range("D22").Value = "'0"
$ws.Range("A22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"
$ws.Range("A22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 2

# Row 23
$ws.Range("F23").Value = "'0"
$ws.Range("A23").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -100
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = 10
$ws.Range("L23").Value = 0

# Row 24
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -21.739130434782
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 14.925373134328
$ws.Range("I24").Value = 305
$ws.Range("J24").Value = 302
$ws.Range("K24").Value = 0.993377483443
$ws.Range("L24").Value = 3.389830508474
$ws.Range("M24").Value = 39.269406392694

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 127.272727272727
$ws.Range("I25").Value = 82
$ws.Range("J25").Value = 69
$ws.Range("K25").Value = 18.840579710144
$ws.Range("L25").Value = 60.78431372549

# Row 26
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 111.111111111111
$ws.Range("F26").Value = 68
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = 83.783783783783
$ws.Range("I26").Value = 232
$ws.Range("J26").Value = 197
$ws.Range("K26").Value = 17.766497461928
$ws.Range("L26").Value = 16.582914572864
$ws.Range("M26").Value = 3.571428571428

# Row 27
$ws.Range("D27").NumberFormat = $fmt_D
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = $fmt_E
$ws.Range("E27").Value = -100
$ws.Range("G27").NumberFormat = $fmt_G
$ws.Range("G27").Value = 2
$ws.Range("H27").NumberFormat = $fmt_H
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = -68.75

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 14
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = -30
$ws.Range("L28").Value = -17.647058823529

# Row 29
$ws.Range("F29").Value = 1
$ws.Range("M29").Value = -46.153846153846
$ws.Range("N29").Value = -89.705882352941

# Row 30
$ws.Range("F30").Value = 1
$ws.Range("M30").Value = -40
$ws.Range("N30").Value = -90.90909090909
